$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Header text updates (rich-text shared strings)
# ---------------------------------------------------------------
# A8: "Volume 30   Number  48" -> "...49"
$ws.Range("A8").Characters(21, 2).Text = "49"

# C9: "Report Covering the Week  11/27/2023  Through  12/3/2023"
#     -> "...12/4/2023  Through  12/10/2023"
# Replace the later (rightmost) substring first so the earlier offset
# doesn't shift.
$ws.Range("C9").Characters(48, 9).Text = "12/10/2023"
$ws.Range("C9").Characters(27, 10).Text = "12/4/2023"

# ---------------------------------------------------------------
# Helper donor cells (row 14 is untouched by this revision, so it is
# a safe source of the "no data" placeholder styling/strings used
# throughout the table: style 14 + shared string "0" (count columns)
# or "***.*" (percent columns)), and of plain numeric styles 15/16.
# ---------------------------------------------------------------
function Set-NoDataNumber($addr) {
    $ws.Range("C14").Copy($ws.Range($addr))
}
function Set-NoDataPercent($addr) {
    $ws.Range("E14").Copy($ws.Range($addr))
}
# When a cell currently rendered as a "no data" placeholder (style 14,
# shared text) needs to become a real number, first paste-special just
# the *format* from an existing numeric donor cell of the right column
# family (style 15 = plain count, style 16 = percent), then set the
# value - this re-uses the existing style record instead of minting a
# new one.
function Set-CountNumber($addr, $value) {
    $ws.Range("I14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $value
}
function Set-PercentNumber($addr, $value) {
    $ws.Range("K14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------
# 2) Row 15 (Burglary)
# ---------------------------------------------------------------
Set-CountNumber "D15" 1
Set-PercentNumber "E15" -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = -37.5

# ---------------------------------------------------------------
# 3) Row 16
# ---------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 85.294117647058
$ws.Range("M16").Value = -35.051546391752
$ws.Range("N16").Value = -83.421052631578

# ---------------------------------------------------------------
# 4) Row 17
# ---------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -44.444444444444
$ws.Range("I17").Value = 96
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = -3.030303030303
$ws.Range("L17").Value = 7.865168539325
$ws.Range("M17").Value = 11.627906976744
$ws.Range("N17").Value = -67.123287671232

# ---------------------------------------------------------------
# 5) Row 18 (becomes a "no data" week-to-date row)
# ---------------------------------------------------------------
Set-NoDataNumber "C18"
Set-NoDataNumber "D18"
Set-NoDataPercent "E18"
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -85.714285714285
$ws.Range("L18").Value = 12.5
$ws.Range("M18").Value = -18.181818181818
$ws.Range("N18").Value = -80.291970802919

# ---------------------------------------------------------------
# 6) Row 19
# ---------------------------------------------------------------
Set-CountNumber "C19" 6
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 17
$ws.Range("H19").Value = -26.086956521739
$ws.Range("I19").Value = 158
$ws.Range("J19").Value = 208
$ws.Range("K19").Value = -24.038461538461
$ws.Range("L19").Value = 30.578512396694
$ws.Range("M19").Value = -27.188940092165
$ws.Range("N19").Value = -29.777777777777

# ---------------------------------------------------------------
# 7) Row 20
# ---------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 800
$ws.Range("L20").Value = 12.280701754386
$ws.Range("N20").Value = -83.376623376623

# ---------------------------------------------------------------
# 8) Row 21 (TOTAL row, bold styles 18/19)
# ---------------------------------------------------------------
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 28.571428571428
$ws.Range("F21").Value = 39
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = -36.065573770491
$ws.Range("I21").Value = 468
$ws.Range("J21").Value = 541
$ws.Range("K21").Value = -13.493530499075
$ws.Range("L21").Value = 23.157894736842
$ws.Range("M21").Value = -16.428571428571
$ws.Range("N21").Value = -72.679509632224

# ---------------------------------------------------------------
# 9) Row 22
# ---------------------------------------------------------------
$ws.Range("M22").Value = -84.615384615384

# ---------------------------------------------------------------
# 10) Row 23
# ---------------------------------------------------------------
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -23.076923076923
$ws.Range("I23").Value = 103
$ws.Range("J23").Value = 110
$ws.Range("K23").Value = -6.363636363636
$ws.Range("L23").Value = 17.045454545454
$ws.Range("M23").Value = 37.333333333333

# ---------------------------------------------------------------
# 11) Row 24
# ---------------------------------------------------------------
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 28
$ws.Range("H24").Value = -36.363636363636
$ws.Range("I24").Value = 565
$ws.Range("J24").Value = 440
$ws.Range("K24").Value = 28.409090909090
$ws.Range("L24").Value = 68.154761904761
$ws.Range("M24").Value = 35.817307692307

# ---------------------------------------------------------------
# 12) Row 25
# ---------------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 600
$ws.Range("F25").Value = 13
$ws.Range("H25").Value = 30
$ws.Range("I25").Value = 151
$ws.Range("J25").Value = 161
$ws.Range("K25").Value = -6.211180124223
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -44.485294117647

# ---------------------------------------------------------------
# 13) Row 26
# ---------------------------------------------------------------
Set-NoDataNumber "C26"
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -75
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -16.666666666666

# ---------------------------------------------------------------
# 14) Row 27
# ---------------------------------------------------------------
Set-NoDataNumber "D27"
Set-NoDataPercent "E27"
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666

# ---------------------------------------------------------------
# 15) Row 28
# ---------------------------------------------------------------
Set-NoDataNumber "C28"
$ws.Range("F28").Value = 1
$ws.Range("L28").Value = -40

# ---------------------------------------------------------------
# 16) Row 29
# ---------------------------------------------------------------
Set-NoDataNumber "C29"
$ws.Range("F29").Value = 1
$ws.Range("L29").Value = -30.769230769230

# ---------------------------------------------------------------
# 17) Row 30
# ---------------------------------------------------------------
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 100
$ws.Range("L30").Value = 0

Write-Host "edits applied"
